$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    3 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    4 = @{ B = 0.6753301551942219; C = 1.667794583268128; D = 0.8054896365839992; E = 0.496779210170732; G = 3.645393585217082 }
    5 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    6 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 8.660232485948974;  G = 13.71653804550039 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
